# foto in webp, curriculum
# - Justify ("both"-aligned) the paragraphs that hold the candidate's name,
#   basic info and activity blurb in the right-hand cell (must run before the
#   table geometry is touched, since resizing the table re-flows/renumbers
#   the document's paragraph collection).
# - Shift the two-column header table slightly left (negative table indent)
#   and re-balance the column widths between the photo cell and the text cell.

$d = $word.ActiveDocument

# --- justify the paragraphs in the second (text) column: name, two spacer
#     lines, basic_info, spacer, activity - all six paragraphs of that cell
#     (document paragraphs 2..7; paragraph 1 is the photo cell, untouched).
for ($i = 2; $i -le 7; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.ParagraphFormat.Alignment = 3
}

$t = $d.Tables.Item(1)

# --- table indent: w:tblInd w:w="-459" w:type="dxa"  (Word works in points: 1pt = 20 dxa)
$t.Rows.LeftIndent = -459 / 20

# --- column / cell widths
#   gridCol 2268 -> 2409 ; gridCol 7087 -> 7053  (same conversion, points = dxa / 20)
$t.Columns.Item(1).Width = 2409 / 20
$t.Columns.Item(2).Width = 7053 / 20
